# Add "Percentage" and "Rank" columns (AL, AM) to the result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - bold/centered header style matches the rest of row 1.
$ws.Range("AL1").Value2 = "Percentage"
$ws.Range("AM1").Value2 = "Rank"
$ws.Range("AK1").Copy()
$ws.Range("AL1:AM1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-13: Percentage = Total Weightage (AI) / Total Marks (AK) * 100
#                 Rank = rank of Percentage among all students, highest first.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 38).Formula = "=AI$r*100/AK$r"
}
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 39).Formula = "=RANK(AL$r,`$AL`$2:`$AL`$13)"
}

# Bake formulas down to plain static values (matches the target workbook,
# which stores Percentage/Rank as literal numbers, not live formulas).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 38).Value2 = $ws.Cells.Item($r, 38).Value2
    $ws.Cells.Item($r, 39).Value2 = $ws.Cells.Item($r, 39).Value2
}
